# "Removed HP & switched to CHEBI-lite"
#
# 1) The three CHEBI-derived rows (CHEBI-BIO-ROLE, CHEBI-CHEM, CHEBI-DRUG-ROLE)
#    get their PURL (column B) switched from the old obolibrary URL to the
#    EBI "chebi_lite" URL.
# 2) The HP row is removed entirely (its whole row is deleted, shifting the
#    rows below it - GSSO - up by one).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newChebiUrl = "https://ftp.ebi.ac.uk/pub/databases/chebi/ontology/chebi_lite.owl"

# Find the header row so we work off column names rather than hard-coded
# positions.
$headerRow = 1
$lastCol = $ws.Cells.Item($headerRow, $ws.Columns.Count()).End(-4159).Column()

$idCol = 1
$purlCol = 2
for ($c = 1; $c -le $lastCol; $c++) {
    $h = $ws.Cells.Item($headerRow, $c).Value()
    if ($h -eq "Ontology ID") { $idCol = $c }
    if ($h -eq "PURL") { $purlCol = $c }
}

$lastRow = $ws.Cells.Item($ws.Rows.Count(), $idCol).End(-4162).Row()

# --- 1) Switch the CHEBI PURLs over to the chebi_lite URL ---
for ($r = $headerRow + 1; $r -le $lastRow; $r++) {
    $ontId = $ws.Cells.Item($r, $idCol).Value()
    if ($ontId -like "CHEBI-*") {
        $ws.Cells.Item($r, $purlCol).Value = $newChebiUrl
    }
}

# --- 2) Remove the HP row entirely ---
for ($r = $lastRow; $r -ge $headerRow + 1; $r--) {
    $ontId = $ws.Cells.Item($r, $idCol).Value()
    if ($ontId -eq "HP") {
        $ws.Rows.Item($r).Delete()
    }
}

# --- 3) Leave the view scrolled/selected near the edited rows, matching
#        where the author was working when they saved.
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
[void]$ws.Range("B26").Select()
